$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Npnt -> Itga8 -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.263451
$ws.Range("H2").Value = 0.790353
$ws.Range("I2").Value = 0.1456293935328523
$ws.Range("J2").Value = 0.1456293935328523
$ws.Range("M2").Value = 0.305927
$ws.Range("N2").Value = 0.917781
$ws.Range("O2").Value = 0.01078151192043695
$ws.Range("P2").Value = 0.01078151192043695
$ws.Range("Q2").Value = 0.080596774077
$ws.Range("R2").Value = 0.7253709666929999
$ws.Range("S2").Value = 0.001570105042340451
$ws.Range("T2").Value = 0.001570105042340451

# Row 3 (ECs -> Npnt -> Itga8 -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.263451
$ws.Range("H3").Value = 0.790353
$ws.Range("I3").Value = 0.1456293935328523
$ws.Range("J3").Value = 0.1456293935328523
$ws.Range("O3").Value = 0.09577486496546363
$ws.Range("P3").Value = 0.09577486496546361
$ws.Range("Q3").Value = 0.7159612873260001
$ws.Range("R3").Value = 6.443651585934001
$ws.Range("S3").Value = 0.01394763550061129
$ws.Range("T3").Value = 0.01394763550061129

# Row 4 (ECs -> Npnt -> Itga8 -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.263451
$ws.Range("H4").Value = 0.790353
$ws.Range("I4").Value = 0.1456293935328523
$ws.Range("J4").Value = 0.1456293935328523
$ws.Range("O4").Value = 0.8934436231140994
$ws.Range("P4").Value = 0.8934436231140994
$ws.Range("Q4").Value = 6.678903142162
$ws.Range("R4").Value = 60.110128279458
$ws.Range("S4").Value = 0.1301116529899005
$ws.Range("T4").Value = 0.1301116529899006

# Row 5 (MuSCs -> Npnt -> Itga8 -> ECs)
$ws.Range("I5").Value = 0.8543706064671477
$ws.Range("J5").Value = 0.8543706064671478
$ws.Range("M5").Value = 0.305927
$ws.Range("N5").Value = 0.917781
$ws.Range("O5").Value = 0.01078151192043695
$ws.Range("P5").Value = 0.01078151192043695
$ws.Range("Q5").Value = 0.4728407712
$ws.Range("R5").Value = 4.2555669408
$ws.Range("S5").Value = 0.009211406878096498
$ws.Range("T5").Value = 0.009211406878096498

# Row 6 (MuSCs -> Npnt -> Itga8 -> FAPs)
$ws.Range("I6").Value = 0.8543706064671477
$ws.Range("J6").Value = 0.8543706064671478
$ws.Range("O6").Value = 0.09577486496546363
$ws.Range("P6").Value = 0.09577486496546361
$ws.Range("S6").Value = 0.08182722946485234
$ws.Range("T6").Value = 0.08182722946485234

# Row 7 (MuSCs -> Npnt -> Itga8 -> MuSCs)
$ws.Range("I7").Value = 0.8543706064671477
$ws.Range("J7").Value = 0.8543706064671478
$ws.Range("O7").Value = 0.8934436231140994
$ws.Range("P7").Value = 0.8934436231140994
$ws.Range("S7").Value = 0.7633319701241988
$ws.Range("T7").Value = 0.763331970124199
